# Slide 8, "Content Placeholder 2" (Shapes.Item(2)):
#   - move the "Listy" paragraph so it lands right after the
#     "mogą być różnych typów" paragraph (i.e. before "Opcje ")
#   - split the run "mogą być różnych typów" into two runs:
#     "mogą być różnych " / "typów"
#   - turn on "Shrink text on overflow" (normAutofit) for the text box

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# NOTE: TextRange.Paragraphs(n, 1) returns the paragraph text *including*
# its trailing paragraph-mark (vbCr), so compare with a "starts with".

# --- 1. Move "Listy" (currently paragraph 1) to just before "Opcje " ---
$listy = $tr.Paragraphs(1, 1)
if (-not $listy.Text.StartsWith("Listy")) {
    throw "Unexpected paragraph 1 text: $($listy.Text)"
}
$listy.Delete() | Out-Null

# After the delete, the paragraphs shift up by one; find "Opcje " again
# (it should now be paragraph 4: Krotki, grupa wartości, mogą być różnych typów, Opcje ...)
$opcje = $tr.Paragraphs(4, 1)
if (-not $opcje.Text.StartsWith("Opcje ")) {
    throw "Unexpected paragraph 4 text: $($opcje.Text)"
}
$opcje.InsertBefore("Listy`r") | Out-Null

# --- 2. Split the "mogą być różnych typów" run in two ---
# It is now paragraph 3 (Krotki, grupa wartości, mogą być różnych typów, Listy, Opcje ...)
$mogaPara = $tr.Paragraphs(3, 1)
if (-not $mogaPara.Text.StartsWith("mogą być różnych typów")) {
    throw "Unexpected paragraph 3 text: $($mogaPara.Text)"
}

$splitWord = "typów"
$typowStart = $mogaPara.Start + ($mogaPara.Text.Length - 1) - $splitWord.Length
$typowRange = $tr.Characters($typowStart, $splitWord.Length)
if ($typowRange.Text -ne $splitWord) {
    throw "Unexpected split range text: $($typowRange.Text)"
}
# Re-assert the same size -> forces the run to split in two while keeping
# identical run formatting on both pieces.
$typowRange.Font.Size = $typowRange.Font.Size

# --- 3. Shrink text on overflow for the placeholder ---
$tf.AutoSize = 2  # ppAutoSizeTextToFitShape -> <a:normAutofit/>
